$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds free-form text (e.g. "24.211.11", "0.9990"),
# not numbers. Force text formatting before writing so Excel does not
# auto-convert/round these into numeric values, then restore the default style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.210.68"
$ws.Range("E2").Value = "  +13.56%  "
$ws.Range("D3").Value = "1.676.82"
$ws.Range("E3").Value = "  +8.43%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "309.32"
$ws.Range("E5").Value = "  +9.61%  "
$ws.Range("D6").Value = "0.9989"
$ws.Range("D7").Value = "0.3734"
$ws.Range("E7").Value = "  +2.73%  "
$ws.Range("D8").Value = "0.3452"
$ws.Range("E8").Value = "  +7.85%  "
$ws.Range("D9").Value = "47.55"
$ws.Range("E9").Value = "  +15.73%  "
$ws.Range("D10").Value = "1.185"
$ws.Range("E10").Value = "  +6.29%  "
$ws.Range("D11").Value = "0.07312"
$ws.Range("E11").Value = "  +5.09%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "20.44"
$ws.Range("E13").Value = "  +7.31%  "
$ws.Range("D14").Value = "6.114"
$ws.Range("E14").Value = "  +6.72%  "
$ws.Range("D15").Value = "6.781"
$ws.Range("E15").Value = "  +5.30%  "
$ws.Range("D16").Value = "1.676.96"
$ws.Range("E16").Value = "  +8.29%  "
$ws.Range("E17").Value = "  +5.05%  "
$ws.Range("D18").Value = "0.9990"
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("D19").Value = "0.06727"
$ws.Range("E19").Value = "  +9.43%  "
$ws.Range("D20").Value = "81.88"
$ws.Range("E20").Value = "  +12.07%  "
$ws.Range("D21").Value = "16.58"
$ws.Range("E21").Value = "  +7.88%  "
$ws.Range("D22").Value = "6.154"
$ws.Range("E22").Value = "  +6.92%  "
$ws.Range("E23").Value = "  +5.03%  "
$ws.Range("D24").Value = "24.177.39"
$ws.Range("E24").Value = "  +13.08%  "
$ws.Range("D25").Value = "2.410"
$ws.Range("E25").Value = "  +3.72%  "
$ws.Range("D26").Value = "3.363"
$ws.Range("E26").Value = "  -9.28%  "
$ws.Range("D27").Value = "2.666"
$ws.Range("E27").Value = "  +15.97%  "
$ws.Range("D28").Value = "151.37"
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("D29").Value = "19.56"
$ws.Range("E29").Value = "  +9.31%  "
$ws.Range("D30").Value = "1.863.14"
$ws.Range("E30").Value = "  +8.33%  "
$ws.Range("D31").Value = "126.81"
$ws.Range("E31").Value = "  +6.44%  "
$ws.Range("D32").Value = "6.447"
$ws.Range("E32").Value = "  +22.47%  "
$ws.Range("D33").Value = "4.092"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "0.9947"
$ws.Range("E34").Value = "  +12.15%  "
$ws.Range("D35").Value = "1.776"
$ws.Range("E35").Value = "  +13.76%  "
$ws.Range("D36").Value = "0.08512"
$ws.Range("E36").Value = "  +5.36%  "
$ws.Range("D37").Value = "12.66"
$ws.Range("E37").Value = "  +16.00%  "
$ws.Range("D38").Value = "0.06493"
$ws.Range("E38").Value = "  +9.88%  "
$ws.Range("D39").Value = "5.376"
$ws.Range("E39").Value = "  +7.22%  "
$ws.Range("D40").Value = "8.908"
$ws.Range("E40").Value = "  +12.21%  "
$ws.Range("D41").Value = "0.02353"
$ws.Range("E41").Value = "  +10.20%  "
$ws.Range("E42").Value = "  +4.30%  "
$ws.Range("D43").Value = "0.2146"
$ws.Range("E43").Value = "  +7.01%  "
$ws.Range("D44").Value = "0.6198"
$ws.Range("E44").Value = "  +12.02%  "
$ws.Range("D45").Value = "0.9986"
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "13.28"
$ws.Range("E46").Value = "  +4.97%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.811"
$ws.Range("E47").Value = "  +6.33%  "
$ws.Range("D48").Value = "0.5951"
$ws.Range("E48").Value = "  +7.92%  "
$ws.Range("D49").Value = "127.43"
$ws.Range("E49").Value = "  +4.14%  "
$ws.Range("D50").Value = "2.033"
$ws.Range("E50").Value = "  +7.66%  "
$ws.Range("D51").Value = "0.07171"
$ws.Range("E51").Value = "  +8.34%  "

$ws.Range("D2:D51").Style = "Normal"
